$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for Adria Avila's extra task, right after her last
# existing task (row 14 "Welcome page animations" / soon "Intro animations").
# This pushes Xavier/Juan/Marti's sections (old rows 15-42) down to 16-43.
$ws.Rows("15:15").Insert()

# New task row: "Update Win/Lose scenes", 1 hour estimation, 1 hour total,
# assigned to Adria Avila (matches the formatting copied down from row 14).
$ws.Range("B15").Value2 = "Update Win/Lose scenes"
$ws.Range("C15").Value2 = "Adria Avila"
$ws.Range("D15").Value2 = "1 hour"
$ws.Range("E15").Value2 = "1 hour"

# Rename the existing "Welcome page animations" task.
$ws.Range("B14").Value2 = "Intro animations"

# Update Adria's totals in the per-person summary table (row 9) to reflect
# the newly added 1 hour task.
$ws.Range("H9").Value2 = "11,5 hours"
$ws.Range("I9").Value2 = "17 hours"

# Match the author's new selection / scroll position.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("F18").Select()
